# Remove the trailing "React.js / ES6" slide (sldId 270, slide27.xml),
# which was the last slide in the deck (index 27 of 27).
$p = $ppt.ActivePresentation
$p.Slides.Item($p.Slides.Count).Delete()
